# The deck's Slide Master currently uses the "Integral" theme
# (ppt/theme/theme2.xml) while the Notes Master uses "Office Theme"
# (ppt/theme/theme1.xml). The commit swaps which theme each master
# uses: the Slide Master ends up with the plain "Office Theme" colour
# palette and the Notes Master ends up with the "Integral" palette.
#
# The font scheme (fontScheme) and format scheme (fmtScheme) blocks are
# byte-identical between the two themes already, so the only real
# content difference is the 12 theme colours (clrScheme) - apply the
# "Office Theme" colour scheme to the presentation's theme via the
# Slide Master's ThemeColorScheme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

# Office Theme colour scheme. ColorFormat.RGB takes a standard COM
# "OLE_COLOR" integer, i.e. 0xBBGGRR (blue in the high byte).
$scheme.Item(1).RGB  = 0        # dk1      000000
$scheme.Item(2).RGB  = 16777215 # lt1      FFFFFF
$scheme.Item(3).RGB  = 6968388  # dk2      44546A
$scheme.Item(4).RGB  = 15132391 # lt2      E7E6E6
$scheme.Item(5).RGB  = 13998939 # accent1  5B9BD5
$scheme.Item(6).RGB  = 3243501  # accent2  ED7D31
$scheme.Item(7).RGB  = 10855845 # accent3  A5A5A5
$scheme.Item(8).RGB  = 49407    # accent4  FFC000
$scheme.Item(9).RGB  = 12874308 # accent5  4472C4
$scheme.Item(10).RGB = 4697456  # accent6  70AD47
$scheme.Item(11).RGB = 12673797 # hlink    0563C1
$scheme.Item(12).RGB = 7491477  # folHlink 954F72

# Best-effort: keep the scheme / theme display names in sync with the
# new palette (no-op on hosts that don't persist these, harmless
# either way).
try { $scheme.Name = "Office" } catch {}
try { $master.Theme.Name = "Office Theme" } catch {}
try { $p.TemplateName = "Office Theme" } catch {}
